$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 1950
$ws.Range("I13").Value = 1950
$ws.Range("K13").Value = 1950
$ws.Range("M13").Value = -1781

$ws.Range("H47").Value = 8605.666999999999
$ws.Range("I47").Value = 8605.666999999999
$ws.Range("K47").Value = 8605.666999999999
$ws.Range("M47").Value = -7633.666999999999

$ws.Range("H51").Value = 16670267
$ws.Range("I51").Value = 4000
$ws.Range("J51").Value = 27781112
$ws.Range("K51").Value = 4000
$ws.Range("L51").Value = 27781112
$ws.Range("M51").Value = -3516
$ws.Range("N51").Value = -27782080

$ws.Range("H64").Value = 3846.4285
$ws.Range("J64").Value = 10000
$ws.Range("L64").Value = 10000
$ws.Range("N64").Value = -10496

$ws.Range("H67").Value = 3846.4285
$ws.Range("J67").Value = 10000
$ws.Range("L67").Value = 10000
$ws.Range("N67").Value = -11716

$ws.Range("H121").Value = 12579.4
$ws.Range("J121").Value = 12579.4
$ws.Range("L121").Value = 37738.2
$ws.Range("N121").Value = -41232.2

$ws.Range("H132").Value = 3109.3333
$ws.Range("I132").Value = 1810.0588
$ws.Range("K132").Value = 5430.1764
$ws.Range("M132").Value = -2900.1764

$ws.Range("H135").Value = 60332.176
$ws.Range("I135").Value = 1915.2727
$ws.Range("J135").Value = 167429.83
$ws.Range("K135").Value = 17237.4543
$ws.Range("L135").Value = 1506868.47
$ws.Range("M135").Value = -14702.4543
$ws.Range("N135").Value = -1511938.47

$ws.Range("H137").Value = 1190.7333
$ws.Range("I137").Value = 1157.3334
$ws.Range("J137").Value = 1240.8334
$ws.Range("K137").Value = 3472.0002
$ws.Range("L137").Value = 3722.5002
$ws.Range("M137").Value = -922.0001999999999
$ws.Range("N137").Value = -8822.5002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 17937.37
$ws.Range("I32").Value = 4925.6875
$ws.Range("K32").Value = 4925.6875
$ws.Range("M32").Value = -4638.6875

$ws.Range("H45").Value = 10404.4
$ws.Range("I45").Value = 16500.285
$ws.Range("J45").Value = 5070.5
$ws.Range("K45").Value = 16500.285
$ws.Range("L45").Value = 5070.5
$ws.Range("M45").Value = -16123.285
$ws.Range("N45").Value = -5824.5

$ws.Range("H61").Value = 2003.6875
$ws.Range("I61").Value = 1912.7273
$ws.Range("K61").Value = 1912.7273
$ws.Range("M61").Value = -1700.7273

$ws.Range("H74").Value = 1801.7646
$ws.Range("I74").Value = 1758.125
$ws.Range("K74").Value = 1758.125
$ws.Range("M74").Value = -884.125

$ws.Range("H77").Value = 1801.7646
$ws.Range("I77").Value = 1758.125
$ws.Range("K77").Value = 8790.625
$ws.Range("M77").Value = -4422.625

$ws.Range("H132").Value = 2006
$ws.Range("I132").Value = 1762.65
$ws.Range("K132").Value = 5287.950000000001
$ws.Range("M132").Value = -2757.950000000001

$ws.Range("H136").Value = 2003.6875
$ws.Range("I136").Value = 1912.7273
$ws.Range("K136").Value = 5738.1819
$ws.Range("M136").Value = -3188.1819

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").Value = ""

$ws.Range("H86").Value = 3607.8667
$ws.Range("I86").Value = 3672.7144
$ws.Range("K86").Value = 3672.7144
$ws.Range("M86").Value = -2549.7144

$ws.Range("H89").Value = 3607.8667
$ws.Range("I89").Value = 3672.7144
$ws.Range("K89").Value = 18363.572
$ws.Range("M89").Value = -12747.572

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10174.208
$ws.Range("I31").Value = 2994.7942
$ws.Range("K31").Value = 2994.7942
$ws.Range("M31").Value = -2699.7942

$ws.Range("H34").Value = 10174.208
$ws.Range("I34").Value = 2994.7942
$ws.Range("K34").Value = 2994.7942
$ws.Range("M34").Value = -2792.7942

$ws.Range("H62").Value = 4250
$ws.Range("J62").Value = 3500
$ws.Range("L62").Value = 3500
$ws.Range("N62").Value = -4748

$ws.Range("H65").Value = 4250
$ws.Range("J65").Value = 3500
$ws.Range("L65").Value = 17500
$ws.Range("N65").Value = -23740

$ws.Range("H99").Value = 29839.8
$ws.Range("J99").Value = 4760
$ws.Range("L99").Value = 4760
$ws.Range("N99").Value = -7756

$ws.Range("H126").Value = 29839.8
$ws.Range("J126").Value = 4760
$ws.Range("L126").Value = 14280
$ws.Range("N126").Value = -19220

$ws.Range("H132").Value = 4512.2085
$ws.Range("I132").Value = 4472.95
$ws.Range("K132").Value = 13418.85
$ws.Range("M132").Value = -10888.85

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1243.0588
$ws.Range("I5").Value = 831.4
$ws.Range("J5").Value = 1831.1428
$ws.Range("K5").Value = 2494.2
$ws.Range("L5").Value = 5493.428400000001
$ws.Range("M5").Value = -2382.2
$ws.Range("N5").Value = -5717.428400000001

$ws.Range("H81").Value = 35723470
$ws.Range("I81").Value = 7365.3335
$ws.Range("K81").Value = 22096.0005
$ws.Range("M81").Value = -20973.0005

$ws.Range("H84").Value = 35723470
$ws.Range("I84").Value = 7365.3335
$ws.Range("K84").Value = 66288.0015
$ws.Range("M84").Value = -60672.0015

$ws.Range("H108").Value = 791.8182
$ws.Range("I108").Value = 791.8182
$ws.Range("K108").Value = 2375.4546
$ws.Range("M108").Value = 504.5454

$ws.Range("H135").Value = 1243.0588
$ws.Range("I135").Value = 831.4
$ws.Range("J135").Value = 1831.1428
$ws.Range("K135").Value = 7482.599999999999
$ws.Range("L135").Value = 16480.2852
$ws.Range("M135").Value = -4947.599999999999
$ws.Range("N135").Value = -21550.2852

$ws.Range("H137").Value = 2647
$ws.Range("I137").Value = 1913.75
$ws.Range("K137").Value = 5741.25
$ws.Range("M137").Value = -641.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7543.0625
$ws.Range("I70").Value = 7565.75
$ws.Range("J70").Value = 7475
$ws.Range("K70").Value = 7565.75
$ws.Range("L70").Value = 7475
$ws.Range("M70").Value = -7295.75
$ws.Range("N70").Value = -8015

$ws.Range("H73").Value = 7543.0625
$ws.Range("I73").Value = 7565.75
$ws.Range("J73").Value = 7475
$ws.Range("K73").Value = 7565.75
$ws.Range("L73").Value = 7475
$ws.Range("M73").Value = -6629.75
$ws.Range("N73").Value = -9347

$ws.Range("H113").Value = 3112.05
$ws.Range("J113").Value = 2665
$ws.Range("L113").Value = 2665
$ws.Range("N113").Value = -7005

$ws.Range("H132").Value = 4826.6
$ws.Range("I132").Value = 4826.6
$ws.Range("K132").Value = 14479.8
$ws.Range("M132").Value = -11949.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 522.75
$ws.Range("I55").Value = 700.8570999999999
$ws.Range("J55").Value = 384.22223
$ws.Range("K55").Value = 700.8570999999999
$ws.Range("L55").Value = 384.22223
$ws.Range("M55").Value = -527.8570999999999
$ws.Range("N55").Value = -730.2222300000001

$ws.Range("H93").Value = 9565.950999999999
$ws.Range("I93").Value = 1486.2632
$ws.Range("K93").Value = 1486.2632
$ws.Range("M93").Value = -238.2632000000001

$ws.Range("H122").Value = 99594.09
$ws.Range("I122").Value = 157367.39
$ws.Range("J122").Value = 5712.5
$ws.Range("K122").Value = 472102.17
$ws.Range("L122").Value = 17137.5
$ws.Range("M122").Value = -469652.17
$ws.Range("N122").Value = -22037.5

$ws.Range("H136").Value = 5307.6665
$ws.Range("I136").Value = 4771
$ws.Range("K136").Value = 14313
$ws.Range("M136").Value = -11763

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1337.5264
$ws.Range("I122").Value = 1300.7222
$ws.Range("K122").Value = 3902.1666
$ws.Range("M122").Value = -1452.1666
